$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_3_3_0"
$ws.Range("B2").Value = 0.4617052126849445
$ws.Range("C2").Value = 0.1148055315503338
$ws.Range("D2").Value = 0.4330422374328121
$ws.Range("E2").Value = 0.3904367761852172
$ws.Range("F2").Value = 0.5957337021827698
$ws.Range("G2").Value = 0.1730260252952576
$ws.Range("H2").Value = 0.6501836180686951
$ws.Range("I2").Value = 0.3975711166858673

$ws.Range("A3").Value = "model_3_3_1"
$ws.Range("B3").Value = 0.4639333408989004
$ws.Range("C3").Value = 0.1189240143350556
$ws.Range("D3").Value = 0.4295566718261252
$ws.Range("E3").Value = 0.3882066065083235
$ws.Range("F3").Value = 0.5932678580284119
$ws.Range("G3").Value = 0.1722210049629211
$ws.Range("H3").Value = 0.6541808843612671
$ws.Range("I3").Value = 0.3990256786346436

$ws.Range("A4").Value = "model_3_3_2"
$ws.Range("B4").Value = 0.465342157210617
$ws.Range("C4").Value = 0.2171985783389684
$ws.Range("D4").Value = 0.4126476500159032
$ws.Range("E4").Value = 0.3898081983557685
$ws.Range("F4").Value = 0.5917087197303772
$ws.Range("G4").Value = 0.1530116200447083
$ws.Range("H4").Value = 0.6735720038414001
$ws.Range("I4").Value = 0.3979810774326324

$ws.Range("A5").Value = "model_3_3_3"
$ws.Range("B5").Value = 0.4659477644836172
$ws.Range("C5").Value = 0.2153753517247847
$ws.Range("D5").Value = 0.412427487907096
$ws.Range("E5").Value = 0.3893364290894554
$ws.Range("F5").Value = 0.5910384058952332
$ws.Range("G5").Value = 0.1533679962158203
$ws.Range("H5").Value = 0.6738245487213135
$ws.Range("I5").Value = 0.3982887864112854

$ws.Range("A6").Value = "model_3_3_4"
$ws.Range("B6").Value = 0.4668214178340007
$ws.Range("C6").Value = 0.1930688155397685
$ws.Range("D6").Value = 0.4153525141898092
$ws.Range("E6").Value = 0.3882175485379242
$ws.Range("F6").Value = 0.5900716185569763
$ws.Range("G6").Value = 0.1577281951904297
$ws.Range("H6").Value = 0.670470118522644
$ws.Range("I6").Value = 0.3990185260772705

$ws.Range("A7").Value = "model_3_3_6"
$ws.Range("B7").Value = 0.4670210191530472
$ws.Range("C7").Value = 0.1978607537258751
$ws.Range("D7").Value = 0.413755435200704
$ws.Range("E7").Value = 0.3876564499568248
$ws.Range("F7").Value = 0.5898507237434387
$ws.Range("G7").Value = 0.1567915081977844
$ws.Range("H7").Value = 0.6723016500473022
$ws.Range("I7").Value = 0.3993844985961914

$ws.Range("A8").Value = "model_3_3_5"
$ws.Range("B8").Value = 0.4670378841984539
$ws.Range("C8").Value = 0.1974575611276226
$ws.Range("D8").Value = 0.4139787708747302
$ws.Range("E8").Value = 0.3877773818301479
$ws.Range("F8").Value = 0.5898320078849792
$ws.Range("G8").Value = 0.1568703204393387
$ws.Range("H8").Value = 0.6720455288887024
$ws.Range("I8").Value = 0.3993056118488312

$ws.Range("A9").Value = "model_3_3_24"
$ws.Range("B9").Value = 0.4671679919239531
$ws.Range("C9").Value = 0.1997256262662129
$ws.Range("D9").Value = 0.3994087759232083
$ws.Range("E9").Value = 0.3760812523336214
$ws.Range("F9").Value = 0.5896880030632019
$ws.Range("G9").Value = 0.1564269959926605
$ws.Range("H9").Value = 0.6887543201446533
$ws.Range("I9").Value = 0.4069340825080872

$ws.Range("A10").Value = "model_3_3_23"
$ws.Range("B10").Value = 0.4674727743280016
$ws.Range("C10").Value = 0.2042905413285254
$ws.Range("D10").Value = 0.400797010590475
$ws.Range("E10").Value = 0.3779545329899361
$ws.Range("F10").Value = 0.5893507599830627
$ws.Range("G10").Value = 0.1555347144603729
$ws.Range("H10").Value = 0.6871622800827026
$ws.Range("I10").Value = 0.4057123064994812

$ws.Range("A11").Value = "model_3_3_7"
$ws.Range("B11").Value = 0.4674878489087108
$ws.Range("C11").Value = 0.2059525678924918
$ws.Range("D11").Value = 0.4131206555702188
$ws.Range("E11").Value = 0.3884148902484051
$ws.Range("F11").Value = 0.5893340706825256
$ws.Range("G11").Value = 0.1552098393440247
$ws.Range("H11").Value = 0.6730296015739441
$ws.Range("I11").Value = 0.3988898396492004

$ws.Range("A12").Value = "model_3_3_22"
$ws.Range("B12").Value = 0.4674994267816381
$ws.Range("C12").Value = 0.2040627518198843
$ws.Range("D12").Value = 0.4012723179401879
$ws.Range("E12").Value = 0.3783114397486236
$ws.Range("F12").Value = 0.5893212556838989
$ws.Range("G12").Value = 0.1555792391300201
$ws.Range("H12").Value = 0.6866171956062317
$ws.Range("I12").Value = 0.4054795205593109

$ws.Range("A13").Value = "model_3_3_20"
$ws.Range("B13").Value = 0.4675603176511379
$ws.Range("C13").Value = 0.2049604510229598
$ws.Range("D13").Value = 0.4018300602817648
$ws.Range("E13").Value = 0.3789155887550122
$ws.Range("F13").Value = 0.589253842830658
$ws.Range("G13").Value = 0.1554037630558014
$ws.Range("H13").Value = 0.685977578163147
$ws.Range("I13").Value = 0.4050854742527008

$ws.Range("A14").Value = "model_3_3_21"
$ws.Range("B14").Value = 0.467561042234916
$ws.Range("C14").Value = 0.2046364946156791
$ws.Range("D14").Value = 0.4017246972942234
$ws.Range("E14").Value = 0.3787770040838243
$ws.Range("F14").Value = 0.5892530083656311
$ws.Range("G14").Value = 0.1554670929908752
$ws.Range("H14").Value = 0.6860984563827515
$ws.Range("I14").Value = 0.4051758944988251

$ws.Range("A15").Value = "model_3_3_19"
$ws.Range("B15").Value = 0.4675936872105221
$ws.Range("C15").Value = 0.2043285453120208
$ws.Range("D15").Value = 0.4023952594382775
$ws.Range("E15").Value = 0.379282871201116
$ws.Range("F15").Value = 0.5892168879508972
$ws.Range("G15").Value = 0.1555272787809372
$ws.Range("H15").Value = 0.6853293776512146
$ws.Range("I15").Value = 0.4048459231853485

$ws.Range("A16").Value = "model_3_3_18"
$ws.Range("B16").Value = 0.4676426793694844
$ws.Range("C16").Value = 0.2041694980051495
$ws.Range("D16").Value = 0.4028760630691227
$ws.Range("E16").Value = 0.3796556303700939
$ws.Range("F16").Value = 0.5891627073287964
$ws.Range("G16").Value = 0.1555583775043488
$ws.Range("H16").Value = 0.6847780346870422
$ws.Range("I16").Value = 0.4046027958393097

$ws.Range("A17").Value = "model_3_3_17"
$ws.Range("B17").Value = 0.4676433913800416
$ws.Range("C17").Value = 0.2034508666022097
$ws.Range("D17").Value = 0.4034504514642131
$ws.Range("E17").Value = 0.3800167793508118
$ws.Range("F17").Value = 0.5891618728637695
$ws.Range("G17").Value = 0.155698835849762
$ws.Range("H17").Value = 0.6841192841529846
$ws.Range("I17").Value = 0.4043672978878021

$ws.Range("A18").Value = "model_3_3_10"
$ws.Range("B18").Value = 0.4677611257640671
$ws.Range("C18").Value = 0.2057716385692395
$ws.Range("D18").Value = 0.4104384728655617
$ws.Range("E18").Value = 0.3861671391884244
$ws.Range("F18").Value = 0.5890316367149353
$ws.Range("G18").Value = 0.1552451997995377
$ws.Range("H18").Value = 0.6761054992675781
$ws.Range("I18").Value = 0.4003558158874512

$ws.Range("A19").Value = "model_3_3_9"
$ws.Range("B19").Value = 0.4678234568387689
$ws.Range("C19").Value = 0.2062973316650339
$ws.Range("D19").Value = 0.4108781446195038
$ws.Range("E19").Value = 0.3866141352949558
$ws.Range("F19").Value = 0.5889626145362854
$ws.Range("G19").Value = 0.1551424413919449
$ws.Range("H19").Value = 0.6756012439727783
$ws.Range("I19").Value = 0.4000642895698547

$ws.Range("A20").Value = "model_3_3_13"
$ws.Range("B20").Value = 0.4678987163869254
$ws.Range("C20").Value = 0.2040574397101128
$ws.Range("D20").Value = 0.4068059974249797
$ws.Range("E20").Value = 0.3828892977114954
$ws.Range("F20").Value = 0.5888792872428894
$ws.Range("G20").Value = 0.1555802673101425
$ws.Range("H20").Value = 0.6802712082862854
$ws.Range("I20").Value = 0.4024937152862549

$ws.Range("A21").Value = "model_3_3_14"
$ws.Range("B21").Value = 0.4679253698600536
$ws.Range("C21").Value = 0.206247519530911
$ws.Range("D21").Value = 0.4060260887177162
$ws.Range("E21").Value = 0.3825913408873871
$ws.Range("F21").Value = 0.5888498425483704
$ws.Range("G21").Value = 0.1551521718502045
$ws.Range("H21").Value = 0.6811656355857849
$ws.Range("I21").Value = 0.402688056230545

$ws.Range("A22").Value = "model_3_3_12"
$ws.Range("B22").Value = 0.4679472646819817
$ws.Range("C22").Value = 0.2077427427229313
$ws.Range("D22").Value = 0.4093986170036547
$ws.Range("E22").Value = 0.3856193559313833
$ws.Range("F22").Value = 0.5888256430625916
$ws.Range("G22").Value = 0.1548599153757095
$ws.Range("H22").Value = 0.6772980093955994
$ws.Range("I22").Value = 0.4007131457328796

$ws.Range("A23").Value = "model_3_3_11"
$ws.Range("B23").Value = 0.467967746654155
$ws.Range("C23").Value = 0.2083109317704459
$ws.Range("D23").Value = 0.4096179560232345
$ws.Range("E23").Value = 0.3858915023002454
$ws.Range("F23").Value = 0.5888028740882874
$ws.Range("G23").Value = 0.1547488570213318
$ws.Range("H23").Value = 0.6770464181900024
$ws.Range("I23").Value = 0.4005356431007385

$ws.Range("A24").Value = "model_3_3_16"
$ws.Range("B24").Value = 0.4680436611572176
$ws.Range("C24").Value = 0.2096715441811153
$ws.Range("D24").Value = 0.4052771703827032
$ws.Range("E24").Value = 0.382515427900557
$ws.Range("F24").Value = 0.5887189507484436
$ws.Range("G24").Value = 0.154482901096344
$ws.Range("H24").Value = 0.6820244193077087
$ws.Range("I24").Value = 0.4027375876903534

$ws.Range("A25").Value = "model_3_3_15"
$ws.Range("B25").Value = 0.4680579473179289
$ws.Range("C25").Value = 0.2094509367923234
$ws.Range("D25").Value = 0.4058455759932309
$ws.Range("E25").Value = 0.3829506528614203
$ws.Range("F25").Value = 0.5887030959129333
$ws.Range("G25").Value = 0.154526025056839
$ws.Range("H25").Value = 0.6813726425170898
$ws.Range("I25").Value = 0.4024536907672882

$ws.Range("A26").Value = "model_3_3_8"
$ws.Range("B26").Value = 0.4681011260690893
$ws.Range("C26").Value = 0.2106822626373275
$ws.Range("D26").Value = 0.4129442105148928
$ws.Range("E26").Value = 0.3890192609098078
$ws.Range("F26").Value = 0.5886553525924683
$ws.Range("G26").Value = 0.154285341501236
$ws.Range("H26").Value = 0.6732319593429565
$ws.Range("I26").Value = 0.3984956443309784
